$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 947.375
$ws.Range("I28").Value = 933
$ws.Range("K28").Value = 933
$ws.Range("M28").Value = -448
# Row 40
$ws.Range("H40").Value = 1882.3
$ws.Range("J40").Value = 3766.6667
$ws.Range("L40").Value = 3766.6667
$ws.Range("N40").Value = -4116.6667
# Row 55
$ws.Range("H55").Value = 169.41667
$ws.Range("I55").Value = 124.666664
$ws.Range("J55").Value = 214.16667
$ws.Range("K55").Value = 124.666664
$ws.Range("L55").Value = 214.16667
$ws.Range("M55").Value = 89.333336
$ws.Range("N55").Value = -642.1666700000001
# Row 74
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = ""
# Row 75
$ws.Range("H75").Value = 38997.5
$ws.Range("I75").Value = 35000
$ws.Range("J75").Value = 42995
$ws.Range("K75").Value = 35000
$ws.Range("L75").Value = 42995
$ws.Range("M75").Value = -34064
$ws.Range("N75").Value = -44867
# Row 77
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = ""
# Row 78
$ws.Range("H78").Value = 38997.5
$ws.Range("I78").Value = 35000
$ws.Range("J78").Value = 42995
$ws.Range("K78").Value = 105000
$ws.Range("L78").Value = 128985
$ws.Range("M78").Value = -100320
$ws.Range("N78").Value = -138345
# Row 86
$ws.Range("H86").Value = 5713
$ws.Range("I86").Value = 3001.5
$ws.Range("K86").Value = 3001.5
$ws.Range("M86").Value = -1878.5
# Row 89
$ws.Range("H89").Value = 5713
$ws.Range("I89").Value = 3001.5
$ws.Range("K89").Value = 15007.5
$ws.Range("M89").Value = -9391.5
# Row 94
$ws.Range("H94").Value = 7831.5
$ws.Range("I94").Value = 1996.3334
$ws.Range("J94").Value = 13666.667
$ws.Range("K94").Value = 1996.3334
$ws.Range("L94").Value = 13666.667
$ws.Range("M94").Value = -1545.3334
$ws.Range("N94").Value = -14568.667
# Row 100
$ws.Range("H100").Value = 1948.6666
$ws.Range("I100").Value = 1924.25
$ws.Range("K100").Value = 1924.25
$ws.Range("M100").Value = -1383.25
# Row 111
$ws.Range("H111").Value = 612.5833
$ws.Range("I111").Value = 643.6
$ws.Range("J111").Value = 590.4286
$ws.Range("K111").Value = 1930.8
$ws.Range("L111").Value = 1771.2858
$ws.Range("M111").Value = 1136.2
$ws.Range("N111").Value = -7905.2858
# Row 113
$ws.Range("H113").Value = 10498.786
$ws.Range("I113").Value = 12855
$ws.Range("K113").Value = 12855
$ws.Range("M113").Value = -9601
# Row 116
$ws.Range("H116").Value = 4804.875
$ws.Range("I116").Value = 5063.3335
$ws.Range("J116").Value = 4649.8
$ws.Range("K116").Value = 5063.3335
$ws.Range("L116").Value = 4649.8
$ws.Range("M116").Value = -1621.3335
$ws.Range("N116").Value = -11533.8
# Row 132
$ws.Range("H132").Value = 1892.4286
$ws.Range("I132").Value = 1807.2307
$ws.Range("K132").Value = 5421.6921
$ws.Range("M132").Value = -2891.6921
# Row 138
$ws.Range("H138").Value = 1873.8334
$ws.Range("I138").Value = 1698.625
$ws.Range("J138").Value = 2224.25
$ws.Range("K138").Value = 5095.875
$ws.Range("L138").Value = 6672.75
$ws.Range("M138").Value = 44.125
$ws.Range("N138").Value = -16952.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3792.6943
$ws.Range("I32").Value = 4016.0625
$ws.Range("K32").Value = 4016.0625
$ws.Range("M32").Value = -3729.0625
# Row 61
$ws.Range("H61").Value = 4881.8887
$ws.Range("I61").Value = 3106
$ws.Range("K61").Value = 3106
$ws.Range("M61").Value = -2894
# Row 110
$ws.Range("H110").Value = 400
$ws.Range("I110").Value = 400
$ws.Range("K110").Value = 400
$ws.Range("M110").Value = 1645
# Row 122
$ws.Range("H122").Value = 400
$ws.Range("I122").Value = 400
$ws.Range("K122").Value = 1200
$ws.Range("M122").Value = 1250
# Row 136
$ws.Range("H136").Value = 4881.8887
$ws.Range("I136").Value = 3106
$ws.Range("K136").Value = 9318
$ws.Range("M136").Value = -6768

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").Value = ""
# Row 32
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = ""

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 593.75
$ws.Range("I16").Value = 625
$ws.Range("K16").Value = 625
$ws.Range("M16").Value = -338
# Row 22
$ws.Range("H22").Value = 2793.625
$ws.Range("I22").Value = 724.75
$ws.Range("J22").Value = 4862.5
$ws.Range("K22").Value = 724.75
$ws.Range("L22").Value = 4862.5
$ws.Range("M22").Value = -374.75
$ws.Range("N22").Value = -5562.5
# Row 29
$ws.Range("H29").Value = 2960.5
$ws.Range("J29").Value = 2960.5
$ws.Range("L29").Value = 2960.5
$ws.Range("N29").Value = -3546.5
# Row 99
$ws.Range("H99").Value = 3936
$ws.Range("I99").Value = 3936
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3936
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2438
$ws.Range("N99").Value = ""
# Row 113
$ws.Range("H113").Value = 593.75
$ws.Range("I113").Value = 625
$ws.Range("K113").Value = 625
$ws.Range("M113").Value = 1545
# Row 126
$ws.Range("H126").Value = 3936
$ws.Range("I126").Value = 3936
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 11808
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -9338
$ws.Range("N126").Value = ""
# Row 132
$ws.Range("H132").Value = 2885.25
$ws.Range("I132").Value = 2885.25
$ws.Range("K132").Value = 8655.75
$ws.Range("M132").Value = -6125.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 48
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = ""
# Row 122
$ws.Range("H122").Value = 26486.428
$ws.Range("I122").Value = 36441.3
$ws.Range("K122").Value = 109323.9
$ws.Range("M122").Value = -106873.9

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 3737.5
$ws.Range("I46").Value = 3265.2778
$ws.Range("K46").Value = 3265.2778
$ws.Range("M46").Value = -3077.2778
# Row 48
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = ""
$ws.Range("N48").Value = ""
# Row 99
$ws.Range("H99").Value = 14500
$ws.Range("I99").Value = 14500
$ws.Range("K99").Value = 14500
$ws.Range("M99").Value = -11505
# Row 132
$ws.Range("H132").Value = 3500
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470
# Row 136
$ws.Range("H136").Value = 2001.3334
$ws.Range("I136").Value = 2001.3334
$ws.Range("K136").Value = 6004.0002
$ws.Range("M136").Value = -3454.0002

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1165.5555
$ws.Range("I96").Value = 997.4
$ws.Range("K96").Value = 997.4
$ws.Range("M96").Value = 375.6
# Row 107
$ws.Range("H107").Value = 491.5
$ws.Range("I107").Value = 253.77777
$ws.Range("K107").Value = 761.33331
$ws.Range("M107").Value = 1158.66669
# Row 122
$ws.Range("H122").Value = 2037.2
$ws.Range("I122").Value = 2117.7368
$ws.Range("K122").Value = 6353.2104
$ws.Range("M122").Value = -3903.2104
# Row 126
$ws.Range("H126").Value = 3066.6667
$ws.Range("J126").Value = 3066.6667
$ws.Range("L126").Value = 9200.000100000001
$ws.Range("N126").Value = -14140.0001
# Row 132
$ws.Range("H132").Value = 2624.875
$ws.Range("I132").Value = 1249.5
$ws.Range("J132").Value = 3083.3333
$ws.Range("K132").Value = 3748.5
$ws.Range("L132").Value = 9249.999899999999
$ws.Range("M132").Value = -1218.5
$ws.Range("N132").Value = -14309.9999
# Row 136
$ws.Range("H136").Value = 1270.875
$ws.Range("I136").Value = 698.4
$ws.Range("K136").Value = 2095.2
$ws.Range("M136").Value = 454.8000000000002
